$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all changed cells to Text format first so numeric-looking strings
# (e.g. '70.079.48', '0.629') are preserved exactly as text, matching the
# original inlineStr cell contents instead of being parsed into floats.
$changedCells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "E8", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "E13", "D14", "E14", "D15", "E15", "B16", "C16", "D16", "E16", "B17", "C17", "D17", "E17", "D18", "E18", "D19", "E19", "E20", "E21", "D22", "E22", "B23", "C23", "D23", "E23", "B24", "C24", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "E28", "D29", "E29", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "E35", "D36", "E36", "D37", "E37", "E38", "D39", "E39", "E40", "D41", "E41", "D42", "E42", "E43", "D44", "E44", "D45", "E45", "E46", "D47", "E47", "D48", "E48", "E49", "E50", "D51", "E51")
foreach ($addr in $changedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '70.079.48'
$ws.Range('E2').Value = '  -1.01%  '
$ws.Range('D3').Value = '3.538.74'
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '616.96'
$ws.Range('E5').Value = '  +5.59%  '
$ws.Range('D6').Value = '185.81'
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('D7').Value = '0.629'
$ws.Range('E7').Value = '  +1.23%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('D10').Value = '0.656'
$ws.Range('E10').Value = '  +1.13%  '
$ws.Range('D11').Value = '53.55'
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('D12').Value = '0.0000307'
$ws.Range('E12').Value = '  -4.24%  '
$ws.Range('E13').Value = '  +0.88%  '
$ws.Range('D14').Value = '4.099.69'
$ws.Range('E14').Value = '  -0.98%  '
$ws.Range('D15').Value = '622.65'
$ws.Range('E15').Value = '  +7.75%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '70.123.97'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').Value = '12.79'
$ws.Range('E17').Value = '  +3.30%  '
$ws.Range('D18').Value = '18.98'
$ws.Range('E18').Value = '  -1.97%  '
$ws.Range('D19').Value = '3.508.97'
$ws.Range('E19').Value = '  -1.35%  '
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('E21').Value = '  -1.27%  '
$ws.Range('D22').Value = '17.59'
$ws.Range('E22').Value = '  -0.93%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Value = '104.03'
$ws.Range('E23').Value = '  +9.23%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').Value = '4.73'
$ws.Range('E24').Value = '  +2.99%  '
$ws.Range('D25').Value = '4.95'
$ws.Range('E25').Value = '  -1.95%  '
$ws.Range('D26').Value = '3.02'
$ws.Range('E26').Value = '  +3.18%  '
$ws.Range('D27').Value = '10.99'
$ws.Range('E27').Value = '  -2.59%  '
$ws.Range('E28').Value = '  +8.81%  '
$ws.Range('D29').Value = '34.30'
$ws.Range('E29').Value = '  +6.45%  '
$ws.Range('E30').Value = '  -3.04%  '
$ws.Range('D31').Value = '12.35'
$ws.Range('E31').Value = '  +0.45%  '
$ws.Range('D32').Value = '0.116'
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('D33').Value = '64.06'
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').Value = '3.62'
$ws.Range('E34').Value = '  +17.17%  '
$ws.Range('E35').Value = '  -2.15%  '
$ws.Range('D36').Value = '530.45'
$ws.Range('E36').Value = '  -4.49%  '
$ws.Range('D37').Value = '0.400'
$ws.Range('E37').Value = '  -2.87%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').Value = '37.31'
$ws.Range('E39').Value = '  -0.68%  '
$ws.Range('E40').Value = '  +5.91%  '
$ws.Range('D41').Value = '0.0₃0778'
$ws.Range('E41').Value = '  -4.20%  '
$ws.Range('D42').Value = '3.525.31'
$ws.Range('E42').Value = '  +2.82%  '
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('D44').Value = '0.0466'
$ws.Range('E44').Value = '  +4.94%  '
$ws.Range('D45').Value = '2.95'
$ws.Range('E45').Value = '  -0.32%  '
$ws.Range('E46').Value = '  +4.46%  '
$ws.Range('D47').Value = '3.35'
$ws.Range('E47').Value = '  -4.15%  '
$ws.Range('D48').Value = '9.07'
$ws.Range('E48').Value = '  -3.06%  '
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('E50').Value = '  -1.78%  '
$ws.Range('D51').Value = '133.89'
$ws.Range('E51').Value = '  -1.35%  '
